$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NCBITaxa")

# Header: "NCBI ID" -> "New"
$ws.Range("B1").Value = "New"

# Clear old NCBI ID numeric values (empty cell allowed for not-new taxa)
$ws.Range("B2").ClearContents()
$ws.Range("B4").ClearContents()

# Row 4 (E coli strain): taxonomic name updates
$ws.Range("E4").Value = "Pseudomonadota"
$ws.Range("H4").Value = "Enterobacterales"

# Row 5 (foraminis record): renamed genus/species/name, new phylum
$ws.Range("A5").Value = "M foraminis"
$ws.Range("E5").Value = "Bacillota"
$ws.Range("J5").Value = "Mesobacillus"

# Row 6 (Unknown E coli strain): mark as new, update taxonomy
$ws.Range("B6").Value = "Yes"
$ws.Range("E6").Value = "Pseudomonadota"
$ws.Range("H6").Value = "Enterobacterales"

# Update active cell selection to match target state
$ws.Range("H12").Select()

$wb.Save()
